$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '42.967.91'
$ws.Range('E2').Value = '  -6.76%  '

$ws.Range('D3').Value = '2.548.18'
$ws.Range('E3').Value = '  -1.87%  '

$ws.Range('E4').Value = '  -0.19%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '297.86'
$ws.Range('E5').Value = '  -4.51%  '

$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '94.02'
$ws.Range('E6').Value = '  -5.05%  '

$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.575'
$ws.Range('E7').Value = '  -4.20%  '

$ws.Range('E8').Value = '  -0.03%  '

$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.550'
$ws.Range('E9').Value = '  -5.90%  '

$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '36.09'
$ws.Range('E10').Value = '  -7.55%  '

$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0806'
$ws.Range('E11').Value = '  -4.22%  '

$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '7.72'
$ws.Range('E12').Value = '  -5.10%  '

$ws.Range('B13').Value = 'TRON'
$ws.Range('C13').Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.107'
$ws.Range('E13').Value = '  +0.32%  '

$ws.Range('B14').Value = 'WrappedliquidstakedEther2.0'
$ws.Range('C14').Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range('D14').Value = '2.937.28'
$ws.Range('E14').Value = '  -1.93%  '

$ws.Range('D15').Value = '2.549.13'
$ws.Range('E15').Value = '  -1.64%  '

$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '0.867'
$ws.Range('E16').Value = '  -5.53%  '

$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '14.14'
$ws.Range('E17').Value = '  -4.87%  '

$ws.Range('D18').Value = '42.913.55'
$ws.Range('E18').Value = '  -7.14%  '

$ws.Range('D19').Value = '0.0₃0977'
$ws.Range('E19').Value = '  -4.01%  '

$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '6.64'
$ws.Range('E20').Value = '  -1.53%  '

$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '12.51'
$ws.Range('E21').Value = '  -2.40%  '

$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '72.30'
$ws.Range('E22').Value = '  -1.04%  '

$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '260.80'
$ws.Range('E23').Value = '  -12.05%  '

$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '2.90'
$ws.Range('E24').Value = '  -5.06%  '

$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '29.51'
$ws.Range('E25').Value = '  -0.12%  '

$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '2.12'
$ws.Range('E26').Value = '  -4.89%  '

$ws.Range('E27').Value = '  +0.06%  '

$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '10.05'
$ws.Range('E28').Value = '  -7.23%  '

$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '2.13'
$ws.Range('E29').Value = '  -3.63%  '

$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '36.37'
$ws.Range('E30').Value = '  -6.28%  '

$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '6.00'
$ws.Range('E31').Value = '  -3.39%  '

$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '151.84'
$ws.Range('E32').Value = '  -2.60%  '

$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '2.16'
$ws.Range('E33').Value = '  -1.46%  '

$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '3.39'
$ws.Range('E34').Value = '  -5.97%  '

$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '2.72'
$ws.Range('E35').Value = '  -2.43%  '

$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.0793'
$ws.Range('E36').Value = '  -5.23%  '

$ws.Range('E37').Value = '  -6.70%  '

$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '24.15'
$ws.Range('E38').Value = '  +14.39%  '

$ws.Range('E39').Value = '  -3.68%  '

$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '16.48'
$ws.Range('E40').Value = '  +4.90%  '

$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '3.46'
$ws.Range('E41').Value = '  -3.68%  '

$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.0311'
$ws.Range('E42').Value = '  -6.50%  '

$ws.Range('B43').Value = 'RenderToken'
$ws.Range('C43').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '3.82'
$ws.Range('E43').Value = '  -3.85%  '

$ws.Range('B44').Value = 'Maker'
$ws.Range('C44').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D44').Value = '2.074.89'
$ws.Range('E44').Value = '  -1.63%  '

$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.998'
$ws.Range('E45').Value = '  -0.10%  '

$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '85.24'
$ws.Range('E46').Value = '  -13.52%  '

$ws.Range('E47').Value = '  +3.22%  '

$ws.Range('D48').Value = '2.791.72'
$ws.Range('E48').Value = '  -2.05%  '

$ws.Range('E49').Value = '  -2.48%  '

$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '104.23'
$ws.Range('E50').Value = '  -4.12%  '

$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '8.67'
$ws.Range('E51').Value = '  -9.01%  '
